# Update the Price (column D) and Volume(1h) (column E) cells for the
# crypto rows (2-51) with refreshed values, as captured by the latest
# GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($ws, $addr, $val)
    # Several Price values look like numbers (e.g. "1.004"), and a plain
    # Value assignment would make Excel coerce them into the Number type.
    # Forcing a Text number format keeps the literal string, and restoring
    # the cell style afterwards avoids leaving a visible formatting change.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "25.955.51"

$ws.Range("D3").Value = "1.641.86"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextCellValue $ws "D5" "215.12"
$ws.Range("E5").Value = "  -0.02%  "

Set-TextCellValue $ws "D6" "0.5071"
$ws.Range("E6").Value = "  +0.99%  "

Set-TextCellValue $ws "D7" "1.004"
$ws.Range("E7").Value = "  +0.12%  "

Set-TextCellValue $ws "D8" "0.2560"
$ws.Range("E8").Value = "  -0.23%  "

Set-TextCellValue $ws "D9" "0.06369"
$ws.Range("E9").Value = "  +0.05%  "

Set-TextCellValue $ws "D10" "19.51"
$ws.Range("E10").Value = "  -0.17%  "

Set-TextCellValue $ws "D11" "0.07782"
$ws.Range("E11").Value = "  +0.59%  "

Set-TextCellValue $ws "D12" "4.283"
$ws.Range("E12").Value = "  +0.91%  "

$ws.Range("D13").Value = "1.641.83"
$ws.Range("E13").Value = "  +0.10%  "

Set-TextCellValue $ws "D14" "0.5454"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").Value = "0.0₅7825"
$ws.Range("E15").Value = "  -0.66%  "

Set-TextCellValue $ws "D16" "64.35"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").Value = "26.002.24"
$ws.Range("E17").Value = "  +0.47%  "

Set-TextCellValue $ws "D18" "1.004"
$ws.Range("E18").Value = "  +0.05%  "

Set-TextCellValue $ws "D19" "197.33"
$ws.Range("E19").Value = "  -2.30%  "

$ws.Range("E20").Value = "  +1.55%  "

Set-TextCellValue $ws "D21" "9.962"
$ws.Range("E21").Value = "  +1.08%  "

Set-TextCellValue $ws "D22" "6.050"
$ws.Range("E22").Value = "  +1.34%  "

Set-TextCellValue $ws "D23" "1.005"
$ws.Range("E23").Value = "  +0.20%  "

Set-TextCellValue $ws "D24" "1.883"
$ws.Range("E24").Value = "  +0.36%  "

Set-TextCellValue $ws "D25" "141.39"
$ws.Range("E25").Value = "  +0.66%  "

Set-TextCellValue $ws "D26" "0.1171"
$ws.Range("E26").Value = "  +3.29%  "

Set-TextCellValue $ws "D27" "6.876"
$ws.Range("E27").Value = "  +1.97%  "

$ws.Range("E28").Value = "  +0.71%  "

Set-TextCellValue $ws "D29" "1.238"
$ws.Range("E29").Value = "  -0.05%  "

Set-TextCellValue $ws "D30" "0.04999"
$ws.Range("E30").Value = "  +0.79%  "

Set-TextCellValue $ws "D31" "3.261"
$ws.Range("E31").Value = "  +0.05%  "

Set-TextCellValue $ws "D32" "3.189"
$ws.Range("E32").Value = "  +0.11%  "

Set-TextCellValue $ws "D33" "1.540"
$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("E34").Value = "  -0.28%  "

Set-TextCellValue $ws "D35" "0.8938"
$ws.Range("E35").Value = "  +0.57%  "

Set-TextCellValue $ws "D36" "2.586"
$ws.Range("E36").Value = "  -1.45%  "

$ws.Range("D37").Value = "1.130.07"
$ws.Range("E37").Value = "  -0.90%  "

Set-TextCellValue $ws "D38" "0.5440"
$ws.Range("E38").Value = "  -3.20%  "

Set-TextCellValue $ws "D39" "0.01556"
$ws.Range("E39").Value = "  -0.25%  "

Set-TextCellValue $ws "D40" "2.553"
$ws.Range("E40").Value = "  -0.61%  "

Set-TextCellValue $ws "D41" "1.004"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  +10.50%  "

Set-TextCellValue $ws "D43" "5.598"
$ws.Range("E43").Value = "  -0.96%  "

Set-TextCellValue $ws "D44" "0.8174"
$ws.Range("E44").Value = "  +1.74%  "

Set-TextCellValue $ws "D45" "99.82"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "1.776.58"
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("E47").Value = "  -0.04%  "

Set-TextCellValue $ws "D48" "1.003"
$ws.Range("E48").Value = "  +0.05%  "

Set-TextCellValue $ws "D49" "54.82"
$ws.Range("E49").Value = "  +0.29%  "

Set-TextCellValue $ws "D50" "0.05074"
$ws.Range("E50").Value = "  +0.42%  "

Set-TextCellValue $ws "D51" "1.002"
$ws.Range("E51").Value = "  +0.11%  "
